# Auto-generated Excel COM-interop script
# Applies the Brynhildr_Profits market-data refresh (scheduled runner update)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 107.7
$ws.Range("I9").Value = 113.46667
$ws.Range("J9").Value = 90.40000000000001
$ws.Range("K9").Value = 113.46667
$ws.Range("L9").Value = 90.40000000000001
$ws.Range("M9").Value = 55.53333000000001
$ws.Range("N9").Value = -428.4
# Row 53
$ws.Range("H53").Value = 242.625
$ws.Range("I53").Value = 233.5
$ws.Range("K53").Value = 233.5
$ws.Range("M53").Value = 403.5
# Row 70
$ws.Range("H70").Value = 4538.4614
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 16500
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -17040
# Row 73
$ws.Range("H73").Value = 4538.4614
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 16500
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -18372
# Row 86
$ws.Range("H86").Value = 5380.5386
$ws.Range("I86").Value = 4601.2593
$ws.Range("K86").Value = 4601.2593
$ws.Range("M86").Value = -3478.2593
# Row 89
$ws.Range("H89").Value = 5380.5386
$ws.Range("I89").Value = 4601.2593
$ws.Range("K89").Value = 23006.2965
$ws.Range("M89").Value = -17390.2965
# Row 116
$ws.Range("H116").Value = 17726
$ws.Range("J116").Value = 14211.714
$ws.Range("L116").Value = 14211.714
$ws.Range("N116").Value = -21095.714

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2943281.8
$ws.Range("I61").Value = 1710.1034
$ws.Range("J61").Value = 20004398
$ws.Range("K61").Value = 1710.1034
$ws.Range("L61").Value = 20004398
$ws.Range("M61").Value = -1498.1034
$ws.Range("N61").Value = -20004822
# Row 97
$ws.Range("H97").Value = 657.1579
$ws.Range("I97").Value = 657.1579
$ws.Range("K97").Value = 657.1579
$ws.Range("M97").Value = -161.1579
# Row 124
$ws.Range("H124").Value = 59714.5
$ws.Range("J124").Value = 59714.5
$ws.Range("L124").Value = 59714.5
$ws.Range("N124").Value = -69534.5
# Row 136
$ws.Range("H136").Value = 2943281.8
$ws.Range("I136").Value = 1710.1034
$ws.Range("J136").Value = 20004398
$ws.Range("K136").Value = 5130.3102
$ws.Range("L136").Value = 60013194
$ws.Range("M136").Value = -2580.3102
$ws.Range("N136").Value = -60018294

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3581.7778
$ws.Range("I86").Value = 1746.1428
$ws.Range("K86").Value = 1746.1428
$ws.Range("M86").Value = -623.1428000000001
# Row 89
$ws.Range("H89").Value = 3581.7778
$ws.Range("I89").Value = 1746.1428
$ws.Range("K89").Value = 8730.714
$ws.Range("M89").Value = -3114.714
# Row 134
$ws.Range("H134").Value = 3088073
$ws.Range("J134").Value = 16668667
$ws.Range("L134").Value = 50006001
$ws.Range("N134").Value = -50011071

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
# Row 124
$ws.Range("H124").Value = 49581.25
$ws.Range("J124").Value = 49581.25
$ws.Range("L124").Value = 49581.25
$ws.Range("N124").Value = -54491.25
# Row 132
$ws.Range("H132").Value = 2435.15
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 834666.5600000001
$ws.Range("I46").Value = 1014.2857
$ws.Range("J46").Value = 2001779.8
$ws.Range("K46").Value = 3042.8571
$ws.Range("L46").Value = 6005339.4
$ws.Range("M46").Value = -2951.8571
$ws.Range("N46").Value = -6005521.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1159.4828
$ws.Range("J97").Value = 1624.625
$ws.Range("L97").Value = 1624.625
$ws.Range("N97").Value = -2616.625
# Row 132
$ws.Range("H132").Value = 19990.5
$ws.Range("I132").Value = 11881.917
$ws.Range("J132").Value = 36207.668
$ws.Range("K132").Value = 35645.751
$ws.Range("L132").Value = 108623.004
$ws.Range("M132").Value = -33115.751
$ws.Range("N132").Value = -113683.004

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 21
$ws.Range("H21").Value = 4007
$ws.Range("J21").Value = 4007
$ws.Range("L21").Value = 4007
$ws.Range("N21").Value = -4355
# Row 22
$ws.Range("H22").Value = 3207.8262
$ws.Range("I22").Value = 2852.8462
$ws.Range("J22").Value = 3669.3
$ws.Range("K22").Value = 2852.8462
$ws.Range("L22").Value = 3669.3
$ws.Range("M22").Value = -2557.8462
$ws.Range("N22").Value = -4259.3
# Row 27
$ws.Range("H27").Value = 3207.8262
$ws.Range("I27").Value = 2852.8462
$ws.Range("J27").Value = 3669.3
$ws.Range("K27").Value = 2852.8462
$ws.Range("L27").Value = 3669.3
$ws.Range("M27").Value = -2745.8462
$ws.Range("N27").Value = -3883.3
# Row 40
$ws.Range("H40").Value = 6793.8
$ws.Range("I40").Value = 6298.1665
$ws.Range("J40").Value = 7537.25
$ws.Range("K40").Value = 6298.1665
$ws.Range("L40").Value = 7537.25
$ws.Range("M40").Value = -6162.1665
$ws.Range("N40").Value = -7809.25
# Row 61
$ws.Range("H61").Value = 8527.056
$ws.Range("I61").Value = 7845.3076
$ws.Range("J61").Value = 10299.6
$ws.Range("K61").Value = 7845.3076
$ws.Range("L61").Value = 10299.6
$ws.Range("M61").Value = -7643.3076
$ws.Range("N61").Value = -10703.6
# Row 68
$ws.Range("H68").Value = 3616.5356
$ws.Range("I68").Value = 2197.9546
$ws.Range("J68").Value = 8818
$ws.Range("K68").Value = 2197.9546
$ws.Range("L68").Value = 8818
$ws.Range("M68").Value = -1448.9546
$ws.Range("N68").Value = -10316
# Row 71
$ws.Range("H71").Value = 3616.5356
$ws.Range("I71").Value = 2197.9546
$ws.Range("J71").Value = 8818
$ws.Range("K71").Value = 10989.773
$ws.Range("L71").Value = 44090
$ws.Range("M71").Value = -7245.773000000001
$ws.Range("N71").Value = -51578
# Row 113
$ws.Range("H113").Value = 8527.056
$ws.Range("I113").Value = 7845.3076
$ws.Range("J113").Value = 10299.6
$ws.Range("K113").Value = 7845.3076
$ws.Range("L113").Value = 10299.6
$ws.Range("M113").Value = -5675.3076
$ws.Range("N113").Value = -14639.6
# Row 132
$ws.Range("H132").Value = 4169988.5
$ws.Range("I132").Value = 8336383
$ws.Range("K132").Value = 25009149
$ws.Range("M132").Value = -25006619
# Row 136
$ws.Range("H136").Value = 22730452
$ws.Range("I136").Value = 13892275
$ws.Range("K136").Value = 41676825
$ws.Range("M136").Value = -41674275

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 39800
$ws.Range("I62").Value = 39800
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 39800
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -39176
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 39800
$ws.Range("I65").Value = 39800
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 199000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -195880
$ws.Range("N65").ClearContents()
# Row 74
$ws.Range("H74").Value = 15625
$ws.Range("I74").Value = 15625
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 15625
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -14689
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 15625
$ws.Range("I77").Value = 15625
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 46875
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -42195
$ws.Range("N77").ClearContents()
# Row 96
$ws.Range("H96").Value = 1761.8572
$ws.Range("I96").Value = 1672.1666
$ws.Range("K96").Value = 1672.1666
$ws.Range("M96").Value = -299.1666
# Row 132
$ws.Range("H132").Value = 83334830
$ws.Range("I132").Value = 83334830
$ws.Range("K132").Value = 250004490
$ws.Range("M132").Value = -250001960
# Row 136
$ws.Range("H136").Value = 15140964
$ws.Range("I136").Value = 7247791.5
$ws.Range("K136").Value = 21743374.5
$ws.Range("M136").Value = -21740824.5
